# The commit re-targets the clock tree for 250 MHz operation. The only
# real input-level change in the workbook is the PL0 clock divisor (used
# twice, for the two PL0-derived rows in the "Zynq Clock Frequency and
# Emulated System Performance" table): it drops from 8 to 6. Every other
# cell touched by the diff is a formula result that recalculates
# automatically from that single input change.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C28").Value = 6
$ws.Range("C30").Value = 6

# Leave the selection where the author left it after editing C30.
$ws.Range("C31").Select()
